$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-33) holds the "Förändrad" date, stored as serial date 45174.
# Update it to 45175 (one day later) for every data row, matching the diff.
$ws.Range("C2:C33").Value = 45175
